# Appends new log-activity rows (536-549) captured on 2024-12-26 to the
# worksheet, matching the freshly exported / unstyled look of the most
# recent existing rows (no explicit cell style, plain string/number/bool types).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 536
$lastNewRow  = 549

# Row 535 already has no explicit cell style (just like all the rows that
# preceded it since row 518); copy that bare formatting down onto the new
# rows so the appended cells do not pick up the column-level styles.
$fmtSrc = $ws.Range("A535:F535")
$fmtSrc.Copy()
$fmtDst = $ws.Range("A" + $firstNewRow + ":F" + $lastNewRow)
$fmtDst.PasteSpecial(-4122)

# Row 536
$ws.Range("A536").Value = "49B6NS"
$ws.Range("B536").Value = "2024-12-26 09:46:40"
$ws.Range("C536").Value = "GET /"
$ws.Range("D536").Value = 401
$ws.Range("E536").Value = $false
$ws.Range("F536").Value = "Eitss... mau ngapain? Akses terbatas!"

# Row 537
$ws.Range("A537").Value = "36KI6F"
$ws.Range("B537").Value = "2024-12-26 09:46:40"
$ws.Range("C537").Value = "GET /favicon.ico"
$ws.Range("D537").Value = 404
$ws.Range("E537").Value = $false
$ws.Range("F537").Value = "Not Found"

# Row 538
$ws.Range("A538").Value = "XKGNVC"
$ws.Range("B538").Value = "2024-12-26 09:47:33"
$ws.Range("C538").Value = "POST /setup"
$ws.Range("D538").Value = 400
$ws.Range("E538").Value = $false
$ws.Range("F538").Value = "Model Embedder untuk 'openai' harus salah satu dari ['text-embedding-3-large', 'text-embedding-3-small']."

# Row 539
$ws.Range("A539").Value = "ANLMHV"
$ws.Range("B539").Value = "2024-12-26 09:47:40"
$ws.Range("C539").Value = "POST /setup"
$ws.Range("D539").Value = 200
$ws.Range("E539").Value = $true
$ws.Range("F539").Value = "Proses penyiapan dokumen berhasil diselesaikan dan embeddings berhasil disimpan pada vector database.
###
llm:openai
###
model_llm:gpt-4o
###
embedder:openai
###
model_embedder:text-embedding-3-large
###
chunk_size:1000
###
chunk_overlap:200
###
total_chunks:1076"

# Row 540
$ws.Range("A540").Value = "L01CA6"
$ws.Range("B540").Value = "2024-12-26 09:50:08"
$ws.Range("C540").Value = "POST /chat"
$ws.Range("D540").Value = 200
$ws.Range("E540").Value = $true
$ws.Range("F540").Value = "OK
###
Question:kapan jadwal pmb 2025?
###
Answer:Salam Harmoni🙏
Jadwal penerimaan mahasiswa baru (PMB) 2025 untuk Universitas Pendidikan Ganesha mengikuti sistem seleksi masuk perguruan tinggi negeri yang diselenggarakan oleh Panitia Seleksi Nasional Penerimaan Mahasiswa Baru (SNPMB). Sistem seleksi ini resmi dibuka pada Rabu, 11 Desember 2024. Untuk informasi lebih lanjut mengenai jadwal dan prosedur pendaftaran, Anda dapat mengunjungi situs resmi Universitas Pendidikan Ganesha."

# Row 541
$ws.Range("A541").Value = "D1YEPJ"
$ws.Range("B541").Value = "2024-12-26 09:50:37"
$ws.Range("C541").Value = "POST /chat"
$ws.Range("D541").Value = 200
$ws.Range("E541").Value = $true
$ws.Range("F541").Value = "OK
###
Question:kapan jadwal penerimaan mahasiswa 2025?
###
Answer:Salam Harmoni🙏
Untuk jadwal penerimaan mahasiswa baru tahun 2025 melalui Seleksi Nasional Penerimaan Mahasiswa Baru (SNPMB), informasi lengkapnya biasanya diumumkan oleh panitia SNPMB. Secara umum, jadwal penerimaan mahasiswa baru mencakup beberapa tahap, seperti pendaftaran, ujian, dan pengumuman hasil.
Untuk informasi yang lebih akurat dan terkini, Anda dapat mengunjungi situs resmi Universitas Pendidikan Ganesha atau situs resmi SNPMB. Jika ada pertanyaan lebih lanjut, silakan tanyakan!"

# Row 542
$ws.Range("A542").Value = "8HO1OY"
$ws.Range("B542").Value = "2024-12-26 09:51:11"
$ws.Range("C542").Value = "POST /chat"
$ws.Range("D542").Value = 200
$ws.Range("E542").Value = $true
$ws.Range("F542").Value = "OK
###
Question:kapan jadwal pmb 2025?
###
Answer:Salam Harmoni🙏
Jadwal penerimaan mahasiswa baru (PMB) 2025 untuk Universitas Pendidikan Ganesha mengikuti sistem seleksi masuk perguruan tinggi negeri yang diselenggarakan oleh Panitia Seleksi Nasional Penerimaan Mahasiswa Baru (SNPMB). Sistem seleksi ini resmi dibuka pada Rabu, 11 Desember 2024. Untuk informasi lebih lanjut mengenai jadwal dan prosedur pendaftaran, Anda dapat mengunjungi situs resmi Universitas Pendidikan Ganesha atau situs SNPMB."

# Row 543
$ws.Range("A543").Value = "PNDI2K"
$ws.Range("B543").Value = "2024-12-26 09:51:28"
$ws.Range("C543").Value = "POST /chat"
$ws.Range("D543").Value = 200
$ws.Range("E543").Value = $true
$ws.Range("F543").Value = "OK
###
Question:kapan jadwal pmb?
###
Answer:Salam Harmoni🙏
Jadwal Penerimaan Mahasiswa Baru (PMB) di Universitas Pendidikan Ganesha dibedakan menjadi 3 jalur, yaitu:
1. Jalur Seleksi Nasional Berbasis Prestasi (SNBP)
2. Jalur Seleksi Nasional Berbasis Tes (SNBT)
3. Jalur Seleksi Mandiri (SMBJM)
Untuk informasi lebih lanjut mengenai waktu spesifik dari masing-masing jalur, silakan merujuk ke sumber resmi yang tersedia."

# Row 544
$ws.Range("A544").Value = "5RJVSP"
$ws.Range("B544").Value = "2024-12-26 09:51:48"
$ws.Range("C544").Value = "POST /chat"
$ws.Range("D544").Value = 200
$ws.Range("E544").Value = $true
$ws.Range("F544").Value = "OK
###
Question:kapan jadwal snbp 2025?
###
Answer:Salam Harmoni🙏
Jadwal Seleksi Nasional Berbasis Prestasi (SNBP) 2025 adalah sebagai berikut:
1. **Pendaftaran SNBP**: 04 – 18 Februari 2025
2. **Pengumuman Hasil SNBP**: 18 Maret 2025
Pastikan untuk mempersiapkan segala sesuatunya sesuai dengan jadwal yang telah ditentukan. Jika ada pertanyaan lebih lanjut, silakan tanyakan!"

# Row 545
$ws.Range("A545").Value = "SN5UVP"
$ws.Range("B545").Value = "2024-12-26 09:52:06"
$ws.Range("C545").Value = "POST /chat"
$ws.Range("D545").Value = 200
$ws.Range("E545").Value = $true
$ws.Range("F545").Value = "OK
###
Question:kapan jadwal snbp 2025?
###
Answer:Salam Harmoni🙏
Jadwal Seleksi Nasional Berbasis Prestasi (SNBP) 2025 adalah sebagai berikut:
1. **Pendaftaran SNBP**: 04 – 18 Februari 2025
2. **Pengumuman Hasil SNBP**: 18 Maret 2025
Jika Anda memerlukan informasi lebih lanjut, silakan tanyakan!"

# Row 546
$ws.Range("A546").Value = "9OHAZ3"
$ws.Range("B546").Value = "2024-12-26 09:52:22"
$ws.Range("C546").Value = "POST /chat"
$ws.Range("D546").Value = 200
$ws.Range("E546").Value = $true
$ws.Range("F546").Value = "OK
###
Question:kapan jadwal snbp?
###
Answer:Salam Harmoni🙏
Jadwal Seleksi Nasional Berbasis Prestasi (SNBP) 2025 adalah sebagai berikut:
1. Pengumuman Kuota Sekolah: 28 Desember 2024
2. Masa Sanggah: 28 Desember 2024 – 17 Januari 2025
3. Registrasi Akun SNPMB Sekolah: 06 – 31 Januari 2025
4. Pengisian PDSS oleh Sekolah: 06 – 31 Januari 2025
5. Registrasi Akun SNPMB Siswa: 13 Januari – 18 Februari 2025
6. Pendaftaran SNBP: 04 – 18 Februari 2025
7. Pengumuman Hasil SNBP: 18 Maret 2025
8. Masa Unduh Kartu Peserta SNBP: 4 Februari – 30 April 2025
Jika ada pertanyaan lebih lanjut, silakan tanyakan!"

# Row 547
$ws.Range("A547").Value = "KN2C1P"
$ws.Range("B547").Value = "2024-12-26 09:52:40"
$ws.Range("C547").Value = "POST /chat"
$ws.Range("D547").Value = 200
$ws.Range("E547").Value = $true
$ws.Range("F547").Value = "OK
###
Question:kapan jadwal snbt?
###
Answer:Salam Harmoni🙏
Jadwal Seleksi Nasional Berbasis Tes (SNBT) adalah sebagai berikut:
1. **Pendaftaran UTBK dan SNBT**: 11 – 27 Maret 2025
2. **Pelaksanaan UTBK**: 
   - Gelombang 1: 23 – 30 April 2025
   - Gelombang 2: 2 – 3 Mei 2025
3. **Pengumuman Hasil Seleksi Jalur SNBT**: 28 Mei 2025
Jika ada pertanyaan lebih lanjut, silakan tanyakan!"

# Row 548
$ws.Range("A548").Value = "WUF4D5"
$ws.Range("B548").Value = "2024-12-26 09:54:05"
$ws.Range("C548").Value = "POST /chat"
$ws.Range("D548").Value = 200
$ws.Range("E548").Value = $true
$ws.Range("F548").Value = "OK
###
Question:rektor undiksha
###
Answer:Salam Harmoni🙏
Rektor Universitas Pendidikan Ganesha (Undiksha) adalah Prof. Dr. I Wayan Lasmawan, M.Pd."

# Row 549
$ws.Range("A549").Value = "NBJ96M"
$ws.Range("B549").Value = "2024-12-26 09:54:20"
$ws.Range("C549").Value = "POST /chat"
$ws.Range("D549").Value = 200
$ws.Range("E549").Value = $true
$ws.Range("F549").Value = "OK
###
Question:siapa rektor undiksha?
###
Answer:Salam Harmoni🙏
Rektor Universitas Pendidikan Ganesha (Undiksha) adalah Prof. Dr. I Wayan Lasmawan, M.Pd."

# Several of the new rows contain multi-line Description text; left alone the
# engine would auto-expand those rows and stamp an explicit ht/customHeight
# on them. Re-autofitting clears that explicit height so the saved XML rows
# look like the source rows (no ht/customHeight attribute at all).
$newRows = $ws.Range($firstNewRow.ToString() + ":" + $lastNewRow.ToString())
$newRows.AutoFit()
